# Slide 3, "Content Placeholder 2" shape: append clarifying text about
# server-side / client-side usage to three existing paragraphs, as
# described in the commit "Added server side and client side references".
#
# Strategy: select the *entire* trailing run of text we want to extend
# (via TextRange.Characters(start,len), matching exactly the existing run's
# known text) and call InsertAfter() on that selection. Because the
# selection's bounds line up exactly with the existing run's bounds, the
# new text is appended into that same run instead of being split into a
# brand-new run -- so paragraphs that originally had multiple runs (e.g.
# the spell-check-flagged "Jquery"/"serializer" runs) keep their run
# layout, and paragraphs with a single run stay single-run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame

function Append-ToParagraphEnd {
    param(
        [int]$ParaIndex,
        [string]$ExistingTailText,
        [string]$NewSuffix
    )
    $tr = $tf.TextRange
    $para = $tr.Paragraphs($ParaIndex)
    $tailLen = $ExistingTailText.Length
    $startPos = $para.Start + $para.Length - 1 - $tailLen
    $sel = $tr.Characters($startPos, $tailLen)
    if ($sel.Text -ne $ExistingTailText) {
        throw "Paragraph $ParaIndex tail mismatch: expected [$ExistingTailText] got [$($sel.Text)]"
    }
    $sel.InsertAfter($NewSuffix) | Out-Null
}

# Paragraph 1: "X3DJSONLD.js ... independent of Jquery." ->
#              "... independent of Jquery.  Server and client side.  Contains some client-side code which shouldn't be used on server."
Append-ToParagraphEnd 1 "." ("  Server and client side.  Contains some client-side code which shouldn" + [char]0x2019 + "t be used on server.")

# Paragraph 2: "loaderJQuery.js ... useful for a web page." ->
#              "... useful for a web page.  Client-side"
Append-ToParagraphEnd 2 "loaderJQuery.js – jQuery and other integrations useful for a web page." "  Client-side"

# Paragraph 3: "convertJSON.js ... Also contains JSON validator." ->
#              "... Also contains JSON validator.  Server-side."
Append-ToParagraphEnd 3 ".  Also contains JSON validator." "  Server-side."
